$wb = $excel.ActiveWorkbook

# ----- Sheet "展览" (Exhibition) -----
$ws1 = $wb.Worksheets.Item("展览")

# Update existing "want-to-go" counts
$ws1.Range("F2").Value = 654
$ws1.Range("F4").Value = 1503

# Append new row 6 with the new event, reusing row 5's cell formatting for column A
$ws1.Range("A5").Copy()
$ws1.Range("A6").PasteSpecial(-4122)
$ws1.Cells.Item(6, 1).Value = 5

$ws1.Cells.Item(6, 2).NumberFormat = "@"
$ws1.Cells.Item(6, 2).Value = "2024-06-09"
$ws1.Cells.Item(6, 2).Style = "Normal"

$ws1.Cells.Item(6, 3).Value = "南宁·布谷鸟动漫展4th"
$ws1.Cells.Item(6, 4).Value = "亭洪路45号 百益上河城"
$ws1.Cells.Item(6, 5).Value = "2024.06.09 10:00-06.10 17:00"
$ws1.Cells.Item(6, 6).Value = 0
$ws1.Cells.Item(6, 7).Value = 35
$ws1.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82241"
$ws1.Cells.Item(6, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg"

# ----- Sheet "全部类型" (All Types) -----
$ws4 = $wb.Worksheets.Item("全部类型")

# Update existing "want-to-go" counts
$ws4.Range("F2").Value = 654
$ws4.Range("F4").Value = 1503

# Append new row 7 with the new event, reusing row 6's cell formatting for column A
$ws4.Range("A6").Copy()
$ws4.Range("A7").PasteSpecial(-4122)
$ws4.Cells.Item(7, 1).Value = 6

$ws4.Cells.Item(7, 2).NumberFormat = "@"
$ws4.Cells.Item(7, 2).Value = "2024-06-09"
$ws4.Cells.Item(7, 2).Style = "Normal"

$ws4.Cells.Item(7, 3).Value = "南宁·布谷鸟动漫展4th"
$ws4.Cells.Item(7, 4).Value = "亭洪路45号 百益上河城"
$ws4.Cells.Item(7, 5).Value = "2024.06.09 10:00-06.10 17:00"
$ws4.Cells.Item(7, 6).Value = 0
$ws4.Cells.Item(7, 7).Value = 35
$ws4.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82241"
$ws4.Cells.Item(7, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg"
